$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "Concept" column header with "Text"
$ws.Range("B1").Value = "Text"

# Replace the concept values ("Variables" / "Experimental Design") with "variables"
$ws.Range("B2").Value = "variables"
$ws.Range("B3").Value = "variables"
$ws.Range("B4").Value = "variables"
$ws.Range("B5").Value = "variables"

# Update the active selection to match the saved view state
$ws.Range("B1:B1048576").Select()
